$d = $word.ActiveDocument

# The document has a single section whose default/first-page headers and
# footers each carry one inline picture:
#   - Headers -> BTec logo ("BTec_Logo-Orange"), originally named
#     "image1.jpg" -> rename to "image2.jpg"
#   - Footers -> Pearson logo (Y:\Together Design\...\PearsonLogo.png),
#     originally named "image2.png" -> rename to "image1.png"
# Walk every header/footer of every section and rename the picture(s) found
# there, identifying them by their (unique, unchanged) description so the
# script is resilient to whichever Headers/Footers index happens to map to
# which physical header/footer part.

foreach ($sec in $d.Sections) {

    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
